# Add a new ticket row (row 9) to the sheet, matching the new support
# request from lorenzo.orozco.garcia@gmail.com ("Problemon" / geocerca).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A9").Value = "lorenzo.orozco.garcia@gmail.com"

# Column B ("Problema") is left blank for this ticket, same as most rows.
$ws.Range("B9").ClearFormats()

$ws.Range("C9").Value = "2024-08-27 16:23:26"
$ws.Range("D9").Value = "pedro"
$ws.Range("E9").Value = "Problemon"
$ws.Range("F9").Value = "No se como hacer geocerca"
